$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above row 31, shifting existing rows down.
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the new error message entry.
$ws.Cells.Item(31, 1).Value = "MSG_ERROR_NOMATCH_MUISCA"
$ws.Cells.Item(31, 2).Value = "Documento: {0} No concuerda con Consignatario y Fecha Aviso Llegada"
